# Appends 6 new data rows (rows 36-41) to each of the 5 worksheets,
# matching the style (date format on column A) used in the existing rows.
$wb = $excel.ActiveWorkbook

# --- Worksheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("A35:G35").Copy()
$ws.Range("A36:G41").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(36, 1).Value = 43938
$ws.Cells.Item(36, 2).Value = 14576
$ws.Cells.Item(36, 3).Value = 1400
$ws.Cells.Item(36, 4).Value = 161
$ws.Cells.Item(36, 5).Value = 336
$ws.Cells.Item(36, 6).Value = 82
$ws.Cells.Item(36, 7).Value = 3459
$ws.Cells.Item(37, 1).Value = 43939
$ws.Cells.Item(37, 2).Value = 15464
$ws.Cells.Item(37, 3).Value = 1511
$ws.Cells.Item(37, 4).Value = 164
$ws.Cells.Item(37, 5).Value = 346
$ws.Cells.Item(37, 6).Value = 90
$ws.Cells.Item(37, 7).Value = 3601
$ws.Cells.Item(38, 1).Value = 43940
$ws.Cells.Item(38, 2).Value = 16060
$ws.Cells.Item(38, 3).Value = 1540
$ws.Cells.Item(38, 4).Value = 165
$ws.Cells.Item(38, 5).Value = 355
$ws.Cells.Item(38, 6).Value = 94
$ws.Cells.Item(38, 7).Value = 3684
$ws.Cells.Item(39, 1).Value = 43941
$ws.Cells.Item(39, 2).Value = 16509
$ws.Cells.Item(39, 3).Value = 1580
$ws.Cells.Item(39, 4).Value = 181
$ws.Cells.Item(39, 5).Value = 364
$ws.Cells.Item(39, 6).Value = 98
$ws.Cells.Item(39, 7).Value = 3751
$ws.Cells.Item(40, 1).Value = 43942
$ws.Cells.Item(40, 2).Value = 17337
$ws.Cells.Item(40, 3).Value = 1765
$ws.Cells.Item(40, 4).Value = 182
$ws.Cells.Item(40, 5).Value = 370
$ws.Cells.Item(40, 6).Value = 141
$ws.Cells.Item(40, 7).Value = 3916
$ws.Cells.Item(41, 1).Value = 43943
$ws.Cells.Item(41, 2).Value = 18100
$ws.Cells.Item(41, 3).Value = 1937
$ws.Cells.Item(41, 4).Value = 187
$ws.Cells.Item(41, 5).Value = 384
$ws.Cells.Item(41, 6).Value = 149
$ws.Cells.Item(41, 7).Value = 4054

# --- Worksheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("A35:G35").Copy()
$ws.Range("A36:G41").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(36, 1).Value = 43938
$ws.Cells.Item(36, 2).Value = 10792.21475571486
$ws.Cells.Item(36, 3).Value = 6974.618665215414
$ws.Cells.Item(36, 4).Value = 1496.599929464595
$ws.Cells.Item(36, 5).Value = 2910.664936207061
$ws.Cells.Item(36, 6).Value = 741.1211427581889
$ws.Cells.Item(36, 7).Value = 10115.73701610499
$ws.Cells.Item(37, 1).Value = 43939
$ws.Cells.Item(37, 2).Value = 11449.69875016291
$ws.Cells.Item(37, 3).Value = 7527.606287957493
$ws.Cells.Item(37, 4).Value = 1524.48688467201
$ws.Cells.Item(37, 5).Value = 2997.291868832272
$ws.Cells.Item(37, 6).Value = 813.425644490695
$ws.Cells.Item(37, 7).Value = 10531.01156258863
$ws.Cells.Item(38, 1).Value = 43940
$ws.Cells.Item(38, 2).Value = 11890.98305274291
$ws.Cells.Item(38, 3).Value = 7672.080531736955
$ws.Cells.Item(38, 4).Value = 1533.782536407814
$ws.Cells.Item(38, 5).Value = 3075.256108194961
$ws.Cells.Item(38, 6).Value = 849.5778953569481
$ws.Cells.Item(38, 7).Value = 10773.74245947695
$ws.Cells.Item(39, 1).Value = 43941
$ws.Cells.Item(39, 2).Value = 12223.42709948522
$ws.Cells.Item(39, 3).Value = 7871.35535074311
$ws.Cells.Item(39, 4).Value = 1682.512964180693
$ws.Cells.Item(39, 5).Value = 3153.22034755765
$ws.Cells.Item(39, 6).Value = 885.7301462232012
$ws.Cells.Item(39, 7).Value = 10969.68185816994
$ws.Cells.Item(40, 1).Value = 43942
$ws.Cells.Item(40, 2).Value = 12836.48649971381
$ws.Cells.Item(40, 3).Value = 8793.001388646575
$ws.Cells.Item(40, 4).Value = 1691.808615916498
$ws.Cells.Item(40, 5).Value = 3205.196507132776
$ws.Cells.Item(40, 6).Value = 1274.366843035422
$ws.Cells.Item(40, 7).Value = 11452.21918330938
$ws.Cells.Item(41, 1).Value = 43943
$ws.Cells.Item(41, 2).Value = 13401.41925620465
$ws.Cells.Item(41, 3).Value = 9649.88311037304
$ws.Cells.Item(41, 4).Value = 1738.286874595523
$ws.Cells.Item(41, 5).Value = 3326.47421280807
$ws.Cells.Item(41, 6).Value = 1346.671344767929
$ws.Cells.Item(41, 7).Value = 11855.79585524418

# --- Worksheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("A35:G35").Copy()
$ws.Range("A36:G41").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(36, 1).Value = 43938
$ws.Cells.Item(36, 2).Value = 847
$ws.Cells.Item(36, 3).Value = 67
$ws.Cells.Item(36, 4).Value = 9
$ws.Cells.Item(36, 5).Value = 15
$ws.Cells.Item(36, 6).Value = 7
$ws.Cells.Item(36, 7).Value = 144
$ws.Cells.Item(37, 1).Value = 43939
$ws.Cells.Item(37, 2).Value = 888
$ws.Cells.Item(37, 3).Value = 111
$ws.Cells.Item(37, 4).Value = 3
$ws.Cells.Item(37, 5).Value = 10
$ws.Cells.Item(37, 6).Value = 8
$ws.Cells.Item(37, 7).Value = 142
$ws.Cells.Item(38, 1).Value = 43940
$ws.Cells.Item(38, 2).Value = 596
$ws.Cells.Item(38, 3).Value = 29
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(38, 5).Value = 9
$ws.Cells.Item(38, 6).Value = 4
$ws.Cells.Item(38, 7).Value = 83
$ws.Cells.Item(39, 1).Value = 43941
$ws.Cells.Item(39, 2).Value = 449
$ws.Cells.Item(39, 3).Value = 40
$ws.Cells.Item(39, 4).Value = 16
$ws.Cells.Item(39, 5).Value = 9
$ws.Cells.Item(39, 6).Value = 4
$ws.Cells.Item(39, 7).Value = 67
$ws.Cells.Item(40, 1).Value = 43942
$ws.Cells.Item(40, 2).Value = 828
$ws.Cells.Item(40, 3).Value = 185
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(40, 5).Value = 6
$ws.Cells.Item(40, 6).Value = 43
$ws.Cells.Item(40, 7).Value = 165
$ws.Cells.Item(41, 1).Value = 43943
$ws.Cells.Item(41, 2).Value = 763
$ws.Cells.Item(41, 3).Value = 172
$ws.Cells.Item(41, 4).Value = 5
$ws.Cells.Item(41, 5).Value = 14
$ws.Cells.Item(41, 6).Value = 8
$ws.Cells.Item(41, 7).Value = 138

# --- Worksheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("A35:G35").Copy()
$ws.Range("A36:G41").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(36, 1).Value = 43938
$ws.Cells.Item(36, 2).Value = 627.1271883980849
$ws.Cells.Item(36, 3).Value = 333.7853218353091
$ws.Cells.Item(36, 4).Value = 83.66086562224443
$ws.Cells.Item(36, 5).Value = 129.9403989378152
$ws.Cells.Item(36, 6).Value = 63.26643901594295
$ws.Cells.Item(36, 7).Value = 421.1234837580569
$ws.Cells.Item(37, 1).Value = 43939
$ws.Cells.Item(37, 2).Value = 657.4839944480512
$ws.Cells.Item(37, 3).Value = 552.9876227420792
$ws.Cells.Item(37, 4).Value = 27.88695520741481
$ws.Cells.Item(37, 5).Value = 86.62693262521016
$ws.Cells.Item(37, 6).Value = 72.30450173250622
$ws.Cells.Item(37, 7).Value = 415.2745464836394
$ws.Cells.Item(38, 1).Value = 43940
$ws.Cells.Item(38, 2).Value = 441.2843025799983
$ws.Cells.Item(38, 3).Value = 144.4742437794621
$ws.Cells.Item(38, 4).Value = 9.295651735804936
$ws.Cells.Item(38, 5).Value = 77.96423936268914
$ws.Cells.Item(38, 6).Value = 36.15225086625311
$ws.Cells.Item(38, 7).Value = 242.7308968883244
$ws.Cells.Item(39, 1).Value = 43941
$ws.Cells.Item(39, 2).Value = 332.4440467423142
$ws.Cells.Item(39, 3).Value = 199.2748190061547
$ws.Cells.Item(39, 4).Value = 148.730427772879
$ws.Cells.Item(39, 5).Value = 77.96423936268914
$ws.Cells.Item(39, 6).Value = 36.15225086625311
$ws.Cells.Item(39, 7).Value = 195.9393986929848
$ws.Cells.Item(40, 1).Value = 43942
$ws.Cells.Item(40, 2).Value = 613.0594002285883
$ws.Cells.Item(40, 3).Value = 921.6460379034654
$ws.Cells.Item(40, 4).Value = 9.295651735804936
$ws.Cells.Item(40, 5).Value = 51.97615957512609
$ws.Cells.Item(40, 6).Value = 388.636696812221
$ws.Cells.Item(40, 7).Value = 482.5373251394402
$ws.Cells.Item(41, 1).Value = 43943
$ws.Cells.Item(41, 2).Value = 564.9327564908368
$ws.Cells.Item(41, 3).Value = 856.8817217264651
$ws.Cells.Item(41, 4).Value = 46.47825867902468
$ws.Cells.Item(41, 5).Value = 121.2777056752942
$ws.Cells.Item(41, 6).Value = 72.30450173250622
$ws.Cells.Item(41, 7).Value = 403.5766719348045

# --- Worksheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("A35:G35").Copy()
$ws.Range("A36:G41").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(36, 1).Value = 43938
$ws.Cells.Item(36, 2).Value = 586.9969716198368
$ws.Cells.Item(36, 3).Value = 499.1834216104174
$ws.Cells.Item(36, 4).Value = 61.35130145631258
$ws.Cells.Item(36, 5).Value = 109.1499351077648
$ws.Cells.Item(36, 6).Value = 46.99792612612904
$ws.Cells.Item(36, 7).Value = 422.2932712129403
$ws.Cells.Item(37, 1).Value = 43939
$ws.Cells.Item(37, 2).Value = 612.3189903249306
$ws.Cells.Item(37, 3).Value = 589.8534642582177
$ws.Cells.Item(37, 4).Value = 55.77391041482961
$ws.Cells.Item(37, 5).Value = 105.6848578027564
$ws.Cells.Item(37, 6).Value = 56.03598884269233
$ws.Cells.Item(37, 7).Value = 455.0473199496781
$ws.Cells.Item(38, 1).Value = 43940
$ws.Cells.Item(38, 2).Value = 585.3680698317899
$ws.Cells.Item(38, 3).Value = 505.161666180602
$ws.Cells.Item(38, 4).Value = 48.33738902618567
$ws.Cells.Item(38, 5).Value = 97.0221645402354
$ws.Cells.Item(38, 6).Value = 54.22837629937967
$ws.Cells.Item(38, 7).Value = 432.2364645794501
$ws.Cells.Item(39, 1).Value = 43941
$ws.Cells.Item(39, 2).Value = 539.1664918435484
$ws.Cells.Item(39, 3).Value = 375.6330338266014
$ws.Cells.Item(39, 4).Value = 57.63304076199059
$ws.Cells.Item(39, 5).Value = 95.28962588773119
$ws.Cells.Item(39, 6).Value = 46.99792612612904
$ws.Cells.Item(39, 7).Value = 360.8794298315571
$ws.Cells.Item(40, 1).Value = 43942
$ws.Cells.Item(40, 2).Value = 534.2797864794076
$ws.Cells.Item(40, 3).Value = 430.4336090532939
$ws.Cells.Item(40, 4).Value = 55.77391041482961
$ws.Cells.Item(40, 5).Value = 84.89439397270598
$ws.Cells.Item(40, 6).Value = 119.3024278586353
$ws.Cells.Item(40, 7).Value = 351.5211301924892
$ws.Cells.Item(41, 1).Value = 43943
$ws.Cells.Item(41, 2).Value = 521.840900097958
$ws.Cells.Item(41, 3).Value = 535.0528890315251
$ws.Cells.Item(41, 4).Value = 48.33738902618565
$ws.Cells.Item(41, 5).Value = 83.16185532020178
$ws.Cells.Item(41, 6).Value = 121.1100404019479
$ws.Cells.Item(41, 7).Value = 348.0117678278387

